$d = $word.ActiveDocument

# --- Locate the sentence boundary that needs to become a paragraph break --
# The target paragraph currently reads (one single paragraph):
#   "segue abaixo a tabela de cálculo ... professor “TCC1”. Lembro que os
#    ajustes indicados ... para o professor de TCC2."
# The edit:
#   1) splits it in two right after "...professor “TCC1”." so that
#      "Lembro que os ajustes..." becomes its own paragraph;
#   2) inserts new text about the "Termo de Compromisso" right before that
#      split point, still inside the first paragraph.

$marker = [char]8221 + ". Lembro"   # closing curly quote + ". Lembro"
$searchRange = $d.Content
$found = $searchRange.Find.Execute($marker, $false, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence boundary."
}

# $searchRange now collapses onto the matched text "”. Lembro" - the split
# point is right after the curly closing quote + period (2 characters in).
$splitPos = $searchRange.Start + 2

# --- Step 1: split the paragraph in two at that point ----------------------
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# The new paragraph now starts with a leading space (the space that used to
# separate the two sentences) - strip it so it starts cleanly with "Lembro".
# Position $splitPos itself now holds the freshly-inserted paragraph mark,
# so the stray space is the character right after it.
$leadingSpace = $d.Range($splitPos + 1, $splitPos + 2)
if ($leadingSpace.Text -eq " ") {
    $leadingSpace.Text = ""
}

# --- Step 2: append the new "Termo de Compromisso" text to the first -------
# paragraph (the one ending in "...“TCC1”."), right before its paragraph
# mark. Each InsertAfter call becomes its own run, matching the diff.
# (Looked up via the split position itself, so we don't depend on a
# hard-coded paragraph index.)
$firstPara = $d.Range($splitPos - 1, $splitPos - 1).Paragraphs(1)

function Append-Run([string]$text) {
    $pos = $firstPara.Range.End - 1
    $rng = $d.Range($pos, $pos)
    $rng.InsertAfter($text)
}

Append-Run " "
Append-Run "E ainda na sequência "
Append-Run "segue "
Append-Run "o Termo"
Append-Run " "
Append-Run ("de Compromisso, as DUAS revisões do seu pré-projeto contendo a avaliação do professor " + [char]8220 + "avaliador" + [char]8221 + " e professor " + [char]8220 + "TCC1" + [char]8221 + ", junto com as avaliações da defesa na banca de qualificação.")
